# feat: add 2022-Q3 data
#
# - "总计" sheet: insert the 2022-Q3 summary as the new row 2, pushing the
#   existing 2021-Q4 summary row down to row 3.
# - Insert a brand-new "2022-Q3" worksheet (positioned between "总计" and
#   "2021-Q4") carrying the per-fund holdings table for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet — shift the existing 2021-Q4 summary row down to row 3 and
#    write the new 2022-Q3 summary into row 2.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Preserve the bold/bordered index-column style (currently on A2) by copying
# its format onto A3 before any values are touched.
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)

# Move the old 2021-Q4 values down to row 3.
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2021-Q4"
$summary.Range("C3").Value = 6
$summary.Range("D3").Value = 0.51

# Write the new 2022-Q3 values into row 2.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.06

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right before the existing "2021-Q4"
#    worksheet.
# ---------------------------------------------------------------------------
$oldQ4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($oldQ4)
$newSheet.Name = "2022-Q3"

# Re-resolve sheets by name now that the tab collection changed shape —
# positional references captured before a structural edit can silently
# repoint to the wrong sheet afterwards.
$q4 = $wb.Worksheets.Item("2021-Q4")
$q3 = $wb.Worksheets.Item("2022-Q3")

# Clone the header/index-column formatting (bold + border on row 1 and the
# A-column index cells) from the 2021-Q4 sheet so the new sheet matches it.
# (Kept as two copies so column A of row 1 — never populated — stays empty
# instead of picking up a stray formatted-but-blank cell.)
$q4.Range("B1:H1").Copy()
$q3.Range("B1").PasteSpecial(-4122)
$q4.Range("A2:A3").Copy()
$q3.Range("A2").PasteSpecial(-4122)

# Headers.
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Row 2 — 014839 兴银碳中和主题混合C.
$q3.Range("A2").Value = 0
$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "014839"
$q3.Range("C2").Value = "兴银碳中和主题混合C"
$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "0.67"
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "92.10"
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "5.06"
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.0339"
$q3.Range("H2").Value = 4

# Row 3 — 014838 兴银碳中和主题混合A.
$q3.Range("A3").Value = 1
$q3.Range("B3").NumberFormat = "@"
$q3.Range("B3").Value = "014838"
$q3.Range("C3").Value = "兴银碳中和主题混合A"
$q3.Range("D3").NumberFormat = "@"
$q3.Range("D3").Value = "0.55"
$q3.Range("E3").NumberFormat = "@"
$q3.Range("E3").Value = "92.10"
$q3.Range("F3").NumberFormat = "@"
$q3.Range("F3").Value = "5.06"
$q3.Range("G3").NumberFormat = "@"
$q3.Range("G3").Value = "0.0278"
$q3.Range("H3").Value = 4

# The text-forcing NumberFormat="@" step above leaves a stray style behind;
# wipe it by re-pasting the clean (General, unstyled, no border) format from
# a never-touched scratch cell onto every cell we touched with "@".
$q3.Range("Z100").Copy()
$q3.Range("B2:B3").PasteSpecial(-4122)
$q3.Range("D2:G3").PasteSpecial(-4122)

Write-Output "done"
